$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.75
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AE3").Value = 19
$ws.Range("AG3").Value = 9

# Row 4
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 3.75
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5

# Row 5
$ws.Range("G5").Value = 1.48
$ws.Range("I5").Value = 7
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("AG5").Value = 19
$ws.Range("BD5").Value = 301
